$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Refresh the time_taken column (F2:F15) with the latest query timestamps ---
$data.Range("F2").Value = "2021-10-05 14:19:33.897156"
$data.Range("F3").Value = "2021-10-05 14:19:33.897164"
$data.Range("F4").Value = "2021-10-05 14:19:33.897167"
$data.Range("F5").Value = "2021-10-05 14:19:33.897170"
$data.Range("F6").Value = "2021-10-05 14:19:33.897173"
$data.Range("F7").Value = "2021-10-05 14:19:33.897176"
$data.Range("F8").Value = "2021-10-05 14:19:33.897178"
$data.Range("F9").Value = "2021-10-05 14:19:33.897181"
$data.Range("F10").Value = "2021-10-05 14:19:33.897184"
$data.Range("F11").Value = "2021-10-05 14:19:33.897187"
$data.Range("F12").Value = "2021-10-05 14:19:33.897189"
$data.Range("F13").Value = "2021-10-05 14:19:33.897192"
$data.Range("F14").Value = "2021-10-05 14:19:33.897194"
$data.Range("F15").Value = "2021-10-05 14:19:33.897197"

# --- Add a new "metadata" sheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (row 1)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Bring over the bold/bordered header style that "data" uses for its own header row
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$data.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

# Data row (row 2)
$meta.Range("A2").Value = 0
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B2").Value = "Choanal atresia"
$meta.Range("C2").Value = 221

# data_version is a text value that merely looks numeric ("1.15") - force text
# storage, then drop back to the default style so no stray formatting lingers.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.15"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2020-10-07T17:36:55.330069Z"
$meta.Range("F2").Value = "2021-10-05 14:19:33.893451"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/221/?format=json"

$meta.Range("A1").Select() | Out-Null
$data.Activate() | Out-Null
